$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free, explicit per-cell updates reflecting the scraped price/volume refresh.
# Cells whose new text could be misread by Excel as a number (and would thus lose
# formatting, e.g. trailing zeros or precision) are forced to text via NumberFormat "@",
# then restored to the default "Normal" style so no stray formatting is introduced.

$ws.Range('D2').Value = '29.770.16'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '1.858.91'
$ws.Range('E3').Value = '  +1.57%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range('D4').Value = '0.9990'
$ws.Range("D4").Style = "Normal"
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '244.54'
$ws.Range("D5").Style = "Normal"
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '0.6430'
$ws.Range("D6").Style = "Normal"
$ws.Range('E6').Value = '  +3.56%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('E8').Value = '  +4.73%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.07539'
$ws.Range("D9").Style = "Normal"
$ws.Range('E10').Value = '  +2.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '24.53'
$ws.Range("D11").Style = "Normal"
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').Value = '1.865.46'
$ws.Range('E13').Value = '  +1.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '5.043'
$ws.Range("D14").Style = "Normal"
$ws.Range('E14').Value = '  +1.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '0.6926'
$ws.Range("D15").Style = "Normal"
$ws.Range('E15').Value = '  +3.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '83.97'
$ws.Range("D16").Style = "Normal"
$ws.Range('E16').Value = '  +1.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '0.000009968'
$ws.Range("D17").Style = "Normal"
$ws.Range('E17').Value = '  +10.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '6.115'
$ws.Range("D18").Style = "Normal"
$ws.Range('E18').Value = '  +4.62%  '
$ws.Range('D19').Value = '29.769.19'
$ws.Range('E19').Value = '  +2.02%  '
$ws.Range('D20').Value = '2.116.96'
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '236.44'
$ws.Range("D21").Style = "Normal"
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('E22').Value = '  +1.74%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '7.520'
$ws.Range("D24").Style = "Normal"
$ws.Range('E24').Value = '  +2.09%  '
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '159.05'
$ws.Range("D26").Style = "Normal"
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '0.1425'
$ws.Range("D27").Style = "Normal"
$ws.Range('E27').Value = '  +2.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '8.539'
$ws.Range("D28").Style = "Normal"
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '17.93'
$ws.Range("D29").Style = "Normal"
$ws.Range('E29').Value = '  +1.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '0.06200'
$ws.Range("D30").Style = "Normal"
$ws.Range('E30').Value = '  +6.50%  '
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '1.289'
$ws.Range("D32").Style = "Normal"
$ws.Range('E32').Value = '  +6.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '4.165'
$ws.Range("D33").Style = "Normal"
$ws.Range('E33').Value = '  +1.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '4.104'
$ws.Range("D34").Style = "Normal"
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '1.888'
$ws.Range("D35").Style = "Normal"
$ws.Range('E35').Value = '  +1.09%  '
$ws.Range('E36').Value = '  +2.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '0.7307'
$ws.Range("D37").Style = "Normal"
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '2.605'
$ws.Range("D38").Style = "Normal"
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '2.825'
$ws.Range("D39").Style = "Normal"
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('E40').Value = '  +2.05%  '
$ws.Range('D41').Value = '1.201.48'
$ws.Range('E41').Value = '  -1.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '0.9226'
$ws.Range("D42").Style = "Normal"
$ws.Range('E42').Value = '  +1.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '6.245'
$ws.Range("D43").Style = "Normal"
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '0.9997'
$ws.Range("D44").Style = "Normal"
$ws.Range('E44').Value = '  -0.19%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '2.024.17'
$ws.Range('E45').Value = '  +1.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '102.07'
$ws.Range("D46").Style = "Normal"
$ws.Range('E46').Value = '  +0.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '66.55'
$ws.Range("D47").Style = "Normal"
$ws.Range('E47').Value = '  +1.61%  '
$ws.Range('E48').Value = '  -1.12%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '0.4065'
$ws.Range("D49").Style = "Normal"
$ws.Range('E49').Value = '  +0.95%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '9.186'
$ws.Range("D50").Style = "Normal"
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '0.05796'
$ws.Range("D51").Style = "Normal"
$ws.Range('E51').Value = '  +0.91%  '

Write-Output "Applied crypto price/volume refresh to $(95) cells."
